$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.223245666666665
$ws.Range("H2").Value = 24.669737
$ws.Range("I2").Value = 0.003010099901484358
$ws.Range("J2").Value = 0.003010099901484359
$ws.Range("M2").Value = 2.839339666666666
$ws.Range("N2").Value = 8.518018999999999
$ws.Range("O2").Value = 0.07557152725297667
$ws.Range("P2").Value = 0.07557152725297665
$ws.Range("Q2").Value = 23.34858761011144
$ws.Range("R2").Value = 210.1372884910029
$ws.Range("S2").Value = 0.0002274778467392075
$ws.Range("T2").Value = 0.0002274778467392075

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.223245666666665
$ws.Range("H3").Value = 24.669737
$ws.Range("I3").Value = 0.003010099901484358
$ws.Range("J3").Value = 0.003010099901484359
$ws.Range("O3").Value = 0.03769827278900544
$ws.Range("P3").Value = 0.03769827278900544
$ws.Range("Q3").Value = 11.64726262600844
$ws.Range("R3").Value = 104.825363634076
$ws.Range("S3").Value = 0.0001134755672083157
$ws.Range("T3").Value = 0.0001134755672083157

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.223245666666665
$ws.Range("H4").Value = 24.669737
$ws.Range("I4").Value = 0.003010099901484358
$ws.Range("J4").Value = 0.003010099901484359
$ws.Range("M4").Value = 32.04516333333333
$ws.Range("N4").Value = 96.13549
$ws.Range("O4").Value = 0.852910260297995
$ws.Range("P4").Value = 0.8529102602979949
$ws.Range("Q4").Value = 263.5152505184589
$ws.Range("R4").Value = 2371.63725466613
$ws.Range("S4").Value = 0.002567345090497993
$ws.Range("T4").Value = 0.002567345090497993

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.223245666666665
$ws.Range("H5").Value = 24.669737
$ws.Range("I5").Value = 0.003010099901484358
$ws.Range("J5").Value = 0.003010099901484359
$ws.Range("M5").Value = 1.270667666666667
$ws.Range("N5").Value = 3.812003
$ws.Range("O5").Value = 0.03381993966002293
$ws.Range("P5").Value = 0.03381993966002293
$ws.Range("Q5").Value = 10.44901238369011
$ws.Range("R5").Value = 94.04111145321099
$ws.Range("S5").Value = 0.0001018013970388419
$ws.Range("T5").Value = 0.000101801397038842

$ws.Range("I6").Value = 0.9908672065823976
$ws.Range("J6").Value = 0.9908672065823977
$ws.Range("M6").Value = 2.839339666666666
$ws.Range("N6").Value = 8.518018999999999
$ws.Range("O6").Value = 0.07557152725297667
$ws.Range("P6").Value = 0.07557152725297665
$ws.Range("Q6").Value = 7685.907624350561
$ws.Range("R6").Value = 69173.16861915504
$ws.Range("S6").Value = 0.07488134810632252
$ws.Range("T6").Value = 0.07488134810632252

$ws.Range("I7").Value = 0.9908672065823976
$ws.Range("J7").Value = 0.9908672065823977
$ws.Range("O7").Value = 0.03769827278900544
$ws.Range("P7").Value = 0.03769827278900544
$ws.Range("S7").Value = 0.03735398225142303
$ws.Range("T7").Value = 0.03735398225142303

$ws.Range("I8").Value = 0.9908672065823976
$ws.Range("J8").Value = 0.9908672065823977
$ws.Range("M8").Value = 32.04516333333333
$ws.Range("N8").Value = 96.13549
$ws.Range("O8").Value = 0.852910260297995
$ws.Range("P8").Value = 0.8529102602979949
$ws.Range("Q8").Value = 86744.17086433797
$ws.Range("R8").Value = 780697.5377790417
$ws.Range("S8").Value = 0.8451208070869399
$ws.Range("T8").Value = 0.8451208070869399

$ws.Range("I9").Value = 0.9908672065823976
$ws.Range("J9").Value = 0.9908672065823977
$ws.Range("M9").Value = 1.270667666666667
$ws.Range("N9").Value = 3.812003
$ws.Range("O9").Value = 0.03381993966002293
$ws.Range("P9").Value = 0.03381993966002293
$ws.Range("Q9").Value = 3439.61464769534
$ws.Range("R9").Value = 30956.53182925806
$ws.Range("S9").Value = 0.03351106913771216
$ws.Range("T9").Value = 0.03351106913771216

$ws.Range("G10").Value = 14.14340733333333
$ws.Range("H10").Value = 42.430222
$ws.Range("I10").Value = 0.005177161275053701
$ws.Range("J10").Value = 0.005177161275053702
$ws.Range("M10").Value = 2.839339666666666
$ws.Range("N10").Value = 8.518018999999999
$ws.Range("O10").Value = 0.07557152725297667
$ws.Range("P10").Value = 0.07557152725297665
$ws.Range("Q10").Value = 40.15793746335755
$ws.Range("R10").Value = 361.421437170218
$ws.Range("S10").Value = 0.0003912459843907762
$ws.Range("T10").Value = 0.0003912459843907762

$ws.Range("G11").Value = 14.14340733333333
$ws.Range("H11").Value = 42.430222
$ws.Range("I11").Value = 0.005177161275053701
$ws.Range("J11").Value = 0.005177161275053702
$ws.Range("O11").Value = 0.03769827278900544
$ws.Range("P11").Value = 0.03769827278900544
$ws.Range("Q11").Value = 20.03247699453955
$ws.Range("R11").Value = 180.292292950856
$ws.Range("S11").Value = 0.0001951700380196496
$ws.Range("T11").Value = 0.0001951700380196496

$ws.Range("G12").Value = 14.14340733333333
$ws.Range("H12").Value = 42.430222
$ws.Range("I12").Value = 0.005177161275053701
$ws.Range("J12").Value = 0.005177161275053702
$ws.Range("M12").Value = 32.04516333333333
$ws.Range("N12").Value = 96.13549
$ws.Range("O12").Value = 0.852910260297995
$ws.Range("P12").Value = 0.8529102602979949
$ws.Range("Q12").Value = 453.2277980865312
$ws.Range("R12").Value = 4079.05018277878
$ws.Range("S12").Value = 0.004415653970710752
$ws.Range("T12").Value = 0.004415653970710752

$ws.Range("G13").Value = 14.14340733333333
$ws.Range("H13").Value = 42.430222
$ws.Range("I13").Value = 0.005177161275053701
$ws.Range("J13").Value = 0.005177161275053702
$ws.Range("M13").Value = 1.270667666666667
$ws.Range("N13").Value = 3.812003
$ws.Range("O13").Value = 0.03381993966002293
$ws.Range("P13").Value = 0.03381993966002293
$ws.Range("Q13").Value = 17.97157039496289
$ws.Range("R13").Value = 161.744133554666
$ws.Range("S13").Value = 0.0001750912819325235
$ws.Range("T13").Value = 0.0001750912819325235

$ws.Range("G14").Value = 2.583085
$ws.Range("H14").Value = 7.749255
$ws.Range("I14").Value = 0.0009455322410643118
$ws.Range("J14").Value = 0.0009455322410643119
$ws.Range("M14").Value = 2.839339666666666
$ws.Range("N14").Value = 8.518018999999999
$ws.Range("O14").Value = 0.07557152725297667
$ws.Range("P14").Value = 0.07557152725297665
$ws.Range("Q14").Value = 7.334255702871666
$ws.Range("R14").Value = 66.008301325845
$ws.Range("S14").Value = 0.00007145531552415974
$ws.Range("T14").Value = 0.00007145531552415974

$ws.Range("G15").Value = 2.583085
$ws.Range("H15").Value = 7.749255
$ws.Range("I15").Value = 0.0009455322410643118
$ws.Range("J15").Value = 0.0009455322410643119
$ws.Range("O15").Value = 0.03769827278900544
$ws.Range("P15").Value = 0.03769827278900544
$ws.Range("Q15").Value = 3.658636820526667
$ws.Range("R15").Value = 32.92773138474
$ws.Range("S15").Value = 0.00003564493235444207
$ws.Range("T15").Value = 0.00003564493235444208

$ws.Range("G16").Value = 2.583085
$ws.Range("H16").Value = 7.749255
$ws.Range("I16").Value = 0.0009455322410643118
$ws.Range("J16").Value = 0.0009455322410643119
$ws.Range("M16").Value = 32.04516333333333
$ws.Range("N16").Value = 96.13549
$ws.Range("O16").Value = 0.852910260297995
$ws.Range("P16").Value = 0.8529102602979949
$ws.Range("Q16").Value = 82.77538072888333
$ws.Range("R16").Value = 744.97842655995
$ws.Range("S16").Value = 0.0008064541498463087
$ws.Range("T16").Value = 0.0008064541498463087

$ws.Range("G17").Value = 2.583085
$ws.Range("H17").Value = 7.749255
$ws.Range("I17").Value = 0.0009455322410643118
$ws.Range("J17").Value = 0.0009455322410643119
$ws.Range("M17").Value = 1.270667666666667
$ws.Range("N17").Value = 3.812003
$ws.Range("O17").Value = 0.03381993966002293
$ws.Range("P17").Value = 0.03381993966002293
$ws.Range("Q17").Value = 3.282242589751666
$ws.Range("R17").Value = 29.540183307765
$ws.Range("S17").Value = 0.00003197784333940127
$ws.Range("T17").Value = 0.00003197784333940128
